$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "95.218.68"
$ws.Range("E2").Value = "  -1.84%  "

# Row 3
$ws.Range("D3").Value = "3.591.95"
$ws.Range("E3").Value = "  -2.79%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "'2.29"
$ws.Range("E5").Value = "  +19.72%  "

# Row 6
$ws.Range("D6").Value = "'224.91"
$ws.Range("E6").Value = "  -5.17%  "

# Row 7
$ws.Range("D7").Value = "'634.02"
$ws.Range("E7").Value = "  -3.78%  "

# Row 8
$ws.Range("D8").Value = "'0.409"
$ws.Range("E8").Value = "  -3.49%  "

# Row 9
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'1.07"
$ws.Range("E9").Value = "  +0.38%  "

# Row 10
$ws.Range("B10").Value = "USDC"
$ws.Range("C10").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D10").Value = "'1.00"
$ws.Range("E10").Value = "  +0.04%  "

# Row 11
$ws.Range("D11").Value = "3.590.70"
$ws.Range("E11").Value = "  -2.80%  "

# Row 12
$ws.Range("D12").Value = "'46.19"
$ws.Range("E12").Value = "  +4.24%  "

# Row 13
$ws.Range("D13").Value = "'0.205"
$ws.Range("E13").Value = "  -1.39%  "

# Row 14
$ws.Range("D14").Value = "'0.0000286"
$ws.Range("E14").Value = "  -4.89%  "

# Row 15
$ws.Range("D15").Value = "'6.44"
$ws.Range("E15").Value = "  -5.19%  "

# Row 16
$ws.Range("D16").Value = "4.261.50"
$ws.Range("E16").Value = "  -2.60%  "

# Row 17
$ws.Range("D17").Value = "94.861.11"
$ws.Range("E17").Value = "  -1.99%  "

# Row 18
$ws.Range("D18").Value = "'8.73"
$ws.Range("E18").Value = "  -4.96%  "

# Row 19
$ws.Range("D19").Value = "'19.95"
$ws.Range("E19").Value = "  +6.21%  "

# Row 20
$ws.Range("D20").Value = "3.589.41"
$ws.Range("E20").Value = "  -3.08%  "

# Row 21
$ws.Range("D21").Value = "'13.01"
$ws.Range("E21").Value = "  -0.15%  "

# Row 22
$ws.Range("D22").Value = "'0.508"
$ws.Range("E22").Value = "  +0.04%  "

# Row 23
$ws.Range("D23").Value = "'499.61"
$ws.Range("E23").Value = "  -3.94%  "

# Row 24
$ws.Range("D24").Value = "'3.22"
$ws.Range("E24").Value = "  -5.91%  "

# Row 25
$ws.Range("D25").Value = "'0.238"
$ws.Range("E25").Value = "  +16.50%  "

# Row 26
$ws.Range("D26").Value = "'116.91"
$ws.Range("E26").Value = "  +15.36%  "

# Row 27
$ws.Range("D27").Value = "'0.0000201"
$ws.Range("E27").Value = "  -4.04%  "

# Row 28
$ws.Range("D28").Value = "'6.71"
$ws.Range("E28").Value = "  -3.57%  "

# Row 29
$ws.Range("D29").Value = "'12.47"
$ws.Range("E29").Value = "  -6.90%  "

# Row 30
$ws.Range("D30").Value = "'12.76"
$ws.Range("E30").Value = "  +1.55%  "

# Row 31
$ws.Range("D31").Value = "'2.88"
$ws.Range("E31").Value = "  -4.67%  "

# Row 32
$ws.Range("E32").Value = "  +0.07%  "

# Row 33
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  -0.11%  "

# Row 34
$ws.Range("D34").Value = "'0.178"
$ws.Range("E34").Value = "  -7.24%  "

# Row 35
$ws.Range("D35").Value = "'1.75"
$ws.Range("E35").Value = "  -6.29%  "

# Row 36
$ws.Range("D36").Value = "'31.46"
$ws.Range("E36").Value = "  -2.29%  "

# Row 37
$ws.Range("D37").Value = "'0.581"
$ws.Range("E37").Value = "  -2.08%  "

# Row 38
$ws.Range("E38").Value = "  -0.03%  "

# Row 39
$ws.Range("D39").Value = "'587.05"
$ws.Range("E39").Value = "  -9.23%  "

# Row 40
$ws.Range("D40").Value = "'8.25"
$ws.Range("E40").Value = "  -6.97%  "

# Row 41
$ws.Range("D41").Value = "'6.79"
$ws.Range("E41").Value = "  -1.50%  "

# Row 42
$ws.Range("D42").Value = "'40.36"
$ws.Range("E42").Value = "  -0.42%  "

# Row 43
$ws.Range("E43").Value = "  -2.40%  "

# Row 44
$ws.Range("D44").Value = "'0.468"
$ws.Range("E44").Value = "  -7.19%  "

# Row 45
$ws.Range("D45").Value = "'1.89"
$ws.Range("E45").Value = "  -8.98%  "

# Row 46
$ws.Range("D46").Value = "'0.0466"
$ws.Range("E46").Value = "  -0.51%  "

# Row 47
$ws.Range("D47").Value = "'0.914"
$ws.Range("E47").Value = "  -5.10%  "

# Row 48
$ws.Range("D48").Value = "'23.44"
$ws.Range("E48").Value = "  -0.77%  "

# Row 49
$ws.Range("D49").Value = "'3.65"
$ws.Range("E49").Value = "  +3.36%  "

# Row 50
$ws.Range("D50").Value = "'222.64"
$ws.Range("E50").Value = "  +8.96%  "

# Row 51
$ws.Range("D51").Value = "'8.46"
$ws.Range("E51").Value = "  -2.57%  "

